# Daily attendance processing: normalise the "Recorded By" column (G) by
# reversing the order of its comma-separated recorder names for every row
# that lists more than one recorder. Single-recorder rows are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Count -gt 1) {
        $n = $parts.Count
        $reversedParts = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $newText = $reversedParts -join ", "
        $cell.Value = $newText
    }
}
